# effort_estimation.xlsx - round numerical values in effort estimation
# calculations for improved accuracy, and extend the scope to cover the
# Shopping Cart + Checkout modules (with their own subfeature breakdown),
# updating the Cost Summary sheet to match the larger effort totals.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Effort Estimation")
$ws2 = $wb.Worksheets.Item("Cost Summary")

# ---------------------------------------------------------------------------
# Sheet "Effort Estimation"
# Rows 2-11 keep the same Module/Feature grouping, but several Subfeature
# names change and every Buffer/Testing figure is rewritten with a clean
# (rounded) decimal instead of the old binary-floating-point noise
# (e.g. 0.8999999999999999 -> 0.9, 0.6000000000000001 -> 0.6).
# ---------------------------------------------------------------------------

# Row 2: User Authentication / User Registration / Form Implementation
$ws1.Range("A2").Value = "User Authentication"
$ws1.Range("B2").Value = "User Registration"
$ws1.Range("C2").Value = "Form Implementation"
$ws1.Range("D2").Value = 5
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = 0.9
$ws1.Range("G2").Value = 7
$ws1.Range("H2").Value = 1.4
$ws1.Range("I2").Value = 1.26

# Row 3: User Authentication / User Registration / Email Verification
$ws1.Range("A3").Value = "User Authentication"
$ws1.Range("B3").Value = "User Registration"
$ws1.Range("C3").Value = "Email Verification"
$ws1.Range("D3").Value = 3
$ws1.Range("E3").Value = 0.6
$ws1.Range("F3").Value = 0.54
$ws1.Range("G3").Value = 4
$ws1.Range("H3").Value = 0.8
$ws1.Range("I3").Value = 0.72

# Row 4: User Authentication / User Registration / Password Encryption
$ws1.Range("A4").Value = "User Authentication"
$ws1.Range("B4").Value = "User Registration"
$ws1.Range("C4").Value = "Password Encryption"
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 0.4
$ws1.Range("F4").Value = 0.36
$ws1.Range("G4").Value = 3
$ws1.Range("H4").Value = 0.6
$ws1.Range("I4").Value = 0.54

# Row 5: User Authentication / User Login / Form Implementation
$ws1.Range("A5").Value = "User Authentication"
$ws1.Range("B5").Value = "User Login"
$ws1.Range("C5").Value = "Form Implementation"
$ws1.Range("D5").Value = 4
$ws1.Range("E5").Value = 0.8
$ws1.Range("F5").Value = 0.72
$ws1.Range("G5").Value = 5
$ws1.Range("H5").Value = 1
$ws1.Range("I5").Value = 0.9

# Row 6: User Authentication / User Login / Session Management
$ws1.Range("A6").Value = "User Authentication"
$ws1.Range("B6").Value = "User Login"
$ws1.Range("C6").Value = "Session Management"
$ws1.Range("D6").Value = 3
$ws1.Range("E6").Value = 0.6
$ws1.Range("F6").Value = 0.54
$ws1.Range("G6").Value = 4
$ws1.Range("H6").Value = 0.8
$ws1.Range("I6").Value = 0.72

# Row 7: Product Catalog / Product Listing / Product Card Design
$ws1.Range("A7").Value = "Product Catalog"
$ws1.Range("B7").Value = "Product Listing"
$ws1.Range("C7").Value = "Product Card Design"
$ws1.Range("D7").Value = 8
$ws1.Range("E7").Value = 1.6
$ws1.Range("F7").Value = 1.44
$ws1.Range("G7").Value = 5
$ws1.Range("H7").Value = 1
$ws1.Range("I7").Value = 0.9

# Row 8: Product Catalog / Product Listing / Filtering and Sorting
$ws1.Range("A8").Value = "Product Catalog"
$ws1.Range("B8").Value = "Product Listing"
$ws1.Range("C8").Value = "Filtering and Sorting"
$ws1.Range("D8").Value = 6
$ws1.Range("E8").Value = 1.2
$ws1.Range("F8").Value = 1.08
$ws1.Range("G8").Value = 4
$ws1.Range("H8").Value = 0.8
$ws1.Range("I8").Value = 0.72

# Row 9: Product Catalog / Product Search / Search Bar Implementation
$ws1.Range("A9").Value = "Product Catalog"
$ws1.Range("B9").Value = "Product Search"
$ws1.Range("C9").Value = "Search Bar Implementation"
$ws1.Range("D9").Value = 5
$ws1.Range("E9").Value = 1
$ws1.Range("F9").Value = 0.9
$ws1.Range("G9").Value = 6
$ws1.Range("H9").Value = 1.2
$ws1.Range("I9").Value = 1.08

# Row 10: Product Catalog / Product Search / Autocomplete Suggestions
$ws1.Range("A10").Value = "Product Catalog"
$ws1.Range("B10").Value = "Product Search"
$ws1.Range("C10").Value = "Autocomplete Suggestions"
$ws1.Range("D10").Value = 4
$ws1.Range("E10").Value = 0.8
$ws1.Range("F10").Value = 0.72
$ws1.Range("G10").Value = 3
$ws1.Range("H10").Value = 0.6
$ws1.Range("I10").Value = 0.54

# Row 11: Shopping Cart / Add to Cart / Add Product to Cart
$ws1.Range("A11").Value = "Shopping Cart"
$ws1.Range("B11").Value = "Add to Cart"
$ws1.Range("C11").Value = "Add Product to Cart"
$ws1.Range("D11").Value = 4
$ws1.Range("E11").Value = 0.8
$ws1.Range("F11").Value = 0.72
$ws1.Range("G11").Value = 3
$ws1.Range("H11").Value = 0.6
$ws1.Range("I11").Value = 0.54

# ---------------------------------------------------------------------------
# Brand-new rows 12-18: the rest of the Shopping Cart module, plus the new
# Checkout module (Shipping Address / Payment).
# ---------------------------------------------------------------------------

# Row 12: Shopping Cart / Add to Cart / Cart Update
$ws1.Range("A12").Value = "Shopping Cart"
$ws1.Range("B12").Value = "Add to Cart"
$ws1.Range("C12").Value = "Cart Update"
$ws1.Range("D12").Value = 5
$ws1.Range("E12").Value = 1
$ws1.Range("F12").Value = 0.9
$ws1.Range("G12").Value = 4
$ws1.Range("H12").Value = 0.8
$ws1.Range("I12").Value = 0.72

# Row 13: Shopping Cart / Cart Summary / Product Quantity Display
$ws1.Range("A13").Value = "Shopping Cart"
$ws1.Range("B13").Value = "Cart Summary"
$ws1.Range("C13").Value = "Product Quantity Display"
$ws1.Range("D13").Value = 3
$ws1.Range("E13").Value = 0.6
$ws1.Range("F13").Value = 0.54
$ws1.Range("G13").Value = 2
$ws1.Range("H13").Value = 0.4
$ws1.Range("I13").Value = 0.36

# Row 14: Shopping Cart / Cart Summary / Total Cost Calculation
$ws1.Range("A14").Value = "Shopping Cart"
$ws1.Range("B14").Value = "Cart Summary"
$ws1.Range("C14").Value = "Total Cost Calculation"
$ws1.Range("D14").Value = 4
$ws1.Range("E14").Value = 0.8
$ws1.Range("F14").Value = 0.72
$ws1.Range("G14").Value = 3
$ws1.Range("H14").Value = 0.6
$ws1.Range("I14").Value = 0.54

# Row 15: Shopping Cart / Remove from Cart / Remove Product from Cart
$ws1.Range("A15").Value = "Shopping Cart"
$ws1.Range("B15").Value = "Remove from Cart"
$ws1.Range("C15").Value = "Remove Product from Cart"
$ws1.Range("D15").Value = 3
$ws1.Range("E15").Value = 0.6
$ws1.Range("F15").Value = 0.54
$ws1.Range("G15").Value = 2
$ws1.Range("H15").Value = 0.4
$ws1.Range("I15").Value = 0.36

# Row 16: Checkout / Shipping Address / Address Form Implementation
$ws1.Range("A16").Value = "Checkout"
$ws1.Range("B16").Value = "Shipping Address"
$ws1.Range("C16").Value = "Address Form Implementation"
$ws1.Range("D16").Value = 6
$ws1.Range("E16").Value = 1.2
$ws1.Range("F16").Value = 1.08
$ws1.Range("G16").Value = 4
$ws1.Range("H16").Value = 0.8
$ws1.Range("I16").Value = 0.72

# Row 17: Checkout / Payment / Payment Gateway Integration
$ws1.Range("A17").Value = "Checkout"
$ws1.Range("B17").Value = "Payment"
$ws1.Range("C17").Value = "Payment Gateway Integration"
$ws1.Range("D17").Value = 8
$ws1.Range("E17").Value = 1.6
$ws1.Range("F17").Value = 1.44
$ws1.Range("G17").Value = 10
$ws1.Range("H17").Value = 2
$ws1.Range("I17").Value = 1.8

# Row 18: Checkout / Payment / Order Confirmation
$ws1.Range("A18").Value = "Checkout"
$ws1.Range("B18").Value = "Payment"
$ws1.Range("C18").Value = "Order Confirmation"
$ws1.Range("D18").Value = 3
$ws1.Range("E18").Value = 0.6
$ws1.Range("F18").Value = 0.54
$ws1.Range("G18").Value = 2
$ws1.Range("H18").Value = 0.4
$ws1.Range("I18").Value = 0.36

# ---------------------------------------------------------------------------
# Row 19: grand Total row (previously row 12, now shifted down because the
# table grew by 7 rows). Column B stays blank, same as before.
# ---------------------------------------------------------------------------
$ws1.Range("A19").Value = "Total"
$ws1.Range("C19").Value = "Total"
$ws1.Range("D19").Value = 76
$ws1.Range("E19").Value = 15.2
$ws1.Range("F19").Value = 13.68
$ws1.Range("G19").Value = 71
$ws1.Range("H19").Value = 14.2
$ws1.Range("I19").Value = 12.78

# ---------------------------------------------------------------------------
# Row 20: Units row (previously row 13).
# ---------------------------------------------------------------------------
$ws1.Range("C20").Value = "Units"
$ws1.Range("D20").Value = "days"
$ws1.Range("E20").Value = "days"
$ws1.Range("F20").Value = "days"
$ws1.Range("G20").Value = "days"
$ws1.Range("H20").Value = "days"
$ws1.Range("I20").Value = "days"

# ---------------------------------------------------------------------------
# Sheet "Cost Summary": the effort-in-days + pricing figures grow along with
# the bigger scope above (Frontend/Backend/Testing day totals and their
# INR pricing, at Rs.15/hr, Rs.16/hr and Rs.12/hr respectively, 8 hrs/day).
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = 100.32
$ws2.Range("D2").Value = "₹12,038.40"

$ws2.Range("B3").Value = 93.72
$ws2.Range("D3").Value = "₹11,996.16"

$ws2.Range("B4").Value = 15.02
$ws2.Range("D4").Value = "₹1,441.92"

$ws2.Range("D5").Value = "₹25,476.48"

# Column B on "Cost Summary" narrows from ~20.7 chars to ~16.7 chars.
$ws2.Columns.Item(2).ColumnWidth = 15.83
